$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
for ($i=1; $i -le 10; $i++) {
    $ws.Columns.Item($i+10).ColumnWidth = $i
}
